# Day-30 SPring boot Examples
# Fills in the next block of the training-log sheet (day-24 .. day-30)
# and extends the sheet with one more trailing blank row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Rows 28-34: new day entries -------------------------------------------
# Column E (dates) needs the same number format as the surrounding rows
# (m/d/yyyy). Copy that formatting from E27 (already using it) before
# writing the new date values so no brand-new style/numFmt gets minted.
$ws.Range("E27").Copy()
$ws.Range("E28:E34").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$rows = @(
    @{ Row = 28; D = "day-24"; E = (Get-Date -Date "2023-05-23"); F = "multi threading ";                                      G = "12:45PM-1:52PM" },
    @{ Row = 29; D = "day-25"; E = (Get-Date -Date "2023-05-24"); F = "Inter Thread Communication,collections(List,Set)";       G = "12:30PM -1:45PM" },
    @{ Row = 30; D = "day-26"; E = (Get-Date -Date "2023-05-25"); F = " Map(Case study)";                                       G = "12:45PM-1:50PM" },
    @{ Row = 31; D = "day-27"; E = (Get-Date -Date "2023-05-26"); F = "Java8 features";                                         G = "1:00pm-2:15pm" },
    @{ Row = 32; D = "day-28"; E = (Get-Date -Date "2023-05-31"); F = "Java8 features";                                         G = "1:00pm-2:15pm" },
    @{ Row = 33; D = "day-29"; E = (Get-Date -Date "2023-06-02"); F = "jpa,hibernate";                                          G = "11:15-2:00Pm" },
    @{ Row = 34; D = "day-30"; E = (Get-Date -Date "2023-06-10"); F = "spring core,spring jpa,springmvc,spring boot case study"; G = "1:05pm-3:45Pm" }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 4).Value = $r.D   # D - Day
    $ws.Cells.Item($n, 5).Value = $r.E   # E - Date
    $ws.Cells.Item($n, 6).Value = $r.F   # F - Concepts
    $ws.Cells.Item($n, 7).Value = $r.G   # G - Timings(IST)
}

# --- Row 49: number-format only tweak (blank date cell becomes date-fmt) --
$ws.Range("E31").Copy()
$ws.Range("E49").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Rows 167/168: extend the trailing blank rows by one ------------------
$ws.Range("D166").Copy()
$ws.Range("D167").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("E167:I167").Copy()
$ws.Range("E168:I168").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- View state: selection moves to G34, viewport scrolled up a bit -------
$ws.Activate()
$excel.Goto($ws.Range("C14"))
$ws.Range("G34").Select()
